$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)
$sa = $sh.SmartArt

$nodes = $sa.AllNodes
for ($i = 1; $i -le $nodes.Count; $i++) {
    $n = $nodes.Item($i)
    if ($n.TextFrame2.TextRange.Text -eq "Rút các bộ ba quan hệ về từ") {
        $n.TextFrame2.TextRange.Text = "Rút trích các bộ ba"
    }
}
